$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit Text
# number format first, otherwise Excel auto-converts the assigned string
# into a numeric value (losing the original textual formatting, e.g.
# trailing zeros such as "555.50" -> 555.5).

$ws.Range('D2').Value = '64.575.19'
$ws.Range('E2').Value = '  +5.43%  '
$ws.Range('D3').Value = '3.083.30'
$ws.Range('E3').Value = '  +3.38%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.50'
$ws.Range('E5').Value = '  +1.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.45'
$ws.Range('E6').Value = '  +9.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.078.86'
$ws.Range('E8').Value = '  +3.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.497'
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.85'
$ws.Range('E10').Value = '  +15.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.152'
$ws.Range('E11').Value = '  +4.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.456'
$ws.Range('E12').Value = '  +2.84%  '
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000226'
$ws.Range('E13').Value = '  +3.61%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.28'
$ws.Range('E14').Value = '  +4.78%  '
$ws.Range('D15').Value = '3.572.50'
$ws.Range('E15').Value = '  +3.20%  '
$ws.Range('D16').Value = '64.445.51'
$ws.Range('E16').Value = '  +5.21%  '
$ws.Range('D17').Value = '3.094.66'
$ws.Range('E17').Value = '  +3.80%  '
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.73'
$ws.Range('E19').Value = '  +2.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '481.97'
$ws.Range('E20').Value = '  +2.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.73'
$ws.Range('E21').Value = '  +5.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.671'
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.47'
$ws.Range('E23').Value = '  +7.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.18'
$ws.Range('E24').Value = '  +10.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.82'
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.80'
$ws.Range('E27').Value = '  +3.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.92'
$ws.Range('E28').Value = '  +4.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.06'
$ws.Range('E29').Value = '  +9.67%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.17'
$ws.Range('E31').Value = '  +3.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.14'
$ws.Range('E32').Value = '  +2.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.45'
$ws.Range('E33').Value = '  +7.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.73'
$ws.Range('E34').Value = '  +4.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '55.29'
$ws.Range('E35').Value = '  +1.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.07'
$ws.Range('E36').Value = '  +4.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '465.14'
$ws.Range('E37').Value = '  +4.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0829'
$ws.Range('E38').Value = '  +5.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0407'
$ws.Range('E39').Value = '  +8.11%  '
$ws.Range('D40').Value = '3.016.92'
$ws.Range('E40').Value = '  -3.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.117'
$ws.Range('E41').Value = '  +1.38%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.27'
$ws.Range('E42').Value = '  +3.00%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.71'
$ws.Range('E43').Value = '  +15.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '27.67'
$ws.Range('E44').Value = '  +8.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.258'
$ws.Range('E45').Value = '  +7.39%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.07'
$ws.Range('E47').Value = '  +6.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.111'
$ws.Range('E48').Value = '  +3.81%  '
$ws.Range('D49').Value = '0.0₃0515'
$ws.Range('E49').Value = '  +7.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '117.20'
$ws.Range('E50').Value = '  +2.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.07'
$ws.Range('E51').Value = '  +4.49%  '
